$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new rows describing the risk and issue tables.
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "risk table"
$ws.Range("C16").Value = "sudung table risk"

$ws.Range("A17").Value = 12
$ws.Range("B17").Value = "issue table"
$ws.Range("C17").Value = "su dung table issue"
$ws.Range("D17").Value = "workUnitID will be projectID"

# Update row 3: change the "Assigment role" note to describe the Type column
# and add the detailed role note (multi-line) in column D.
$ws.Range("C3").Value = "Xai bang Assigment role la Type"
$ws.Range("D3").Value = "0: ProjectOwer + PM, 1 : pM, 2: dev, 3: test, 4 QA, 5 cus, 6: project Owner`nProject Owner chi co' quyen read only va change PM,`nUser tao project se~ mang role la 0;"

# Give the note enough room to read: wrap the text and grow the row.
$ws.Range("D3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

# Move the active selection to the new note cell.
$ws.Range("D3").Select()
